$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 'b"+\x8d\xffM\xe1Ti\xd1S^\xc4s\xc4\xec\xaegF\x03hB|NWLm\xcb[\xc9\xcb\x9c\x9ei\xa5\xbb|\x9d\x84[''\xb91\x83\x90\x1e\xd2K\xf2[\x035\x9e\xeb\xf7\xa3\x93\xf5\xb0&<\xe6v\x96\x89\x8e\xb8?%\x82\\\xca\xe7oV\x0c\xadU\xd2\x9b\x00\x1b#\xe07\x8c\xed\xfc\xb2\xb5{\xcb)\xd49\xaf\x99\xac^\xc3\xdbO\x0bR\xbfZ&\xa5\x8b^\xe5\x7f\x0e\xc5\x04\xac\xa4\x97\x17\xeaC\x1a\xf2\xc0[7\x81Z\x02\xe7\xae\xf7e3[\xf9py\xde\xdf4\xd3\xea\x868R\xe4\tSPqN\xa4\xd6.Z\xe37\xf8\n\xfa\x84\x94R\xfa\x9dX\x06F\xa5SF}h\xab\x93\xc0$\xb7\xe7\x1b(\xe6\x98.\x14/\x7f\xca\x83\xf6\x8c\xd7J\xc9\x0e\xaf\x03B&\xeaB=f\x98\xe2\x92L\x91\xa0n\x89:\xf1g\xeb%\xaeUK(\xba\xc6\xbef\x16zehG\xea=\xaaS\xc4o\xae|\x08\xf5\xbe]T\x8dd\xef\xa4L-N\xe8D\n\x0c\xe5\xb4y\xa8"'
$ws.Range("C4").Value = 'b''gAAAAABlOtRCmpIkP4Z9uUZnV8kFi6TSCh4UnjqS7_t3qfft3lLKmjOLzrauGuXgD2gnN_sjIa5YzEBsrt2wSQDJBNpBYu03Qg=='''
$ws.Range("B5").Value = 'b''\x15H\xd3\xd3\xcf\xec|\xf6\xac\xf9\xd3\xe9\xe9\x94\xf1\xeaa-T\x9f\xc2.e\x15\xafd\xd1\xce=L\xc9\xc3\xe8\x1e4-\xd3R\x9a\xc6\xa5X-\x16a2w\x0fz\xcbFw\x9f\x1f\xcePf\xbb\xf8\x88\xdeA\xfd\xba\x14J\x9dC\xe3\x7f>\x94F\xcf\x0bh\xa7$\xcd\xaeNI)\xa2\xe4<\xfe&\xce\x8a\x94`N8;\xd3\x18\xe4\x1c`\x82J>\xd2\xae\x8f\xae\xd4\xce\xc7\x96\xcdeb\x15VA\x99\xc93\x0f/\xeb\xf6\x9f\x19|\xb0\xf6\x98\xb1\xc9\xb3\x99u\xb3<\x17*\x92%\xb5\x13s\xe5\xf4\xb0$n\xa9Q\x9b\xeb\x1f\x8f\t_\x0f\xaf{\xa0i\xfd\xda\x8a\xa5$\xc4\x89\x90\x99b\x1a\xdc\xe5\xca{f\x12\xde\x06\xc0orJ\xeaH\xa5H\xfb\xbb\xb8\x8c5zV\xa4N\xc29\xf0\xa8\x9a\xaemC\x853\xda\x15\x01\xa6\xaa5I\xb0+\xa6\xa7\x02<\x9c\x11\x1a\xe0\x7f_\xf0\xa9\x1d0v\xa2s3f\x9ds\x83\xb1l\xba\xf0O\x80D\xa6\xf8\xd0F\xf1\xf3\x04L\x9d '''
$ws.Range("C5").Value = 'b''gAAAAABlOtPiomb0JFqbBEbxYATf_cEGJXa7OTvG5QaKI1o_XQuA91Yi7lw1xc2cGFEzMGJL4p1rOYLoOmD_zg9aLZK6YIUmybt5dSNE1plcxb2QJBXkBdw='''
